$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Investor 6) - fill in advisor notification fields:
# H8 = Phone, I8 = WhatsApp Enabled, J8 = Approved, K8 = Send Confirmation Email
$ws.Range("H8").Value = 1234567789
$ws.Range("I8").Value = "Yes"
$ws.Range("J8").Value = "Yes"
$ws.Range("K8").Value = "No"

# Move the selection to I9, matching where the cursor ended up after entry
$ws.Range("I9").Select()
